$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 13: update B13/C13 (both cells already had values, order is safe)
$ws.Range("B13").Value = 4715
$ws.Range("C13").Value = 4724

# Row 14: B14 had same style as column default -> fully clear (cell disappears)
$ws.Range("B14").Clear()

# Row 15: B15 has its own style (distinct from column default) -> keep empty cell w/ style
$ws.Range("B15").ClearContents()

# Row 16: new split times added. Set C before B to avoid calc-engine ordering quirk
# (writing the second-of-two-previously-blank operands last produces a stale result).
$ws.Range("C16").Value = 6028
$ws.Range("B16").Value = 6019

# Row 17: remove current-run marker, style matches column default -> fully clear
$ws.Range("B17").Clear()

# Row 18
$ws.Range("B18").Clear()

# Row 19
$ws.Range("B19").Clear()

# Row 21
$ws.Range("B21").Clear()

# Row 22
$ws.Range("B22").Clear()

# Row 23
$ws.Range("B23").Clear()

# Row 24
$ws.Range("B24").Clear()

# Row 25
$ws.Range("B25").Clear()

# Row 26
$ws.Range("B26").Clear()

$excel.Calculate()

# Move the active selection to B17 (where the runner's current split marker now sits)
$null = $ws.Range("B17").Select()
